$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting existing rows 10:86 down to 11:87
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly record
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value = "Arica y Parinacota"
$ws.Range("D10").Value = 44817
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 100112040
$ws.Range("G10").Value = "Cilantro"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 700
$ws.Range("L10").Value = 800
$ws.Range("M10").Value = 750
$ws.Range("N10").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O10").Value = "Región de Arica y Parinacota"
$ws.Range("P10").Value = 375
$ws.Range("Q10").Value = 2
$ws.Range("R10").Value = "Hortaliza"
